# Generate Report for Archive
# Two files that were previously handed off ("Ready for handoff") have now
# moved back into translation, so flip their Status to "In Translation" on
# every sheet. The third file in that batch (de573e3a...) stays untouched.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de) hold the status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = $newStatus
$overview.Range("C8").Value = $newStatus
$overview.Range("B9").Value = $newStatus
$overview.Range("C9").Value = $newStatus

# --- zh-cn sheet: column C holds the status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C8").Value = $newStatus
$zhcn.Range("C9").Value = $newStatus

# --- de-de sheet: column C holds the status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C8").Value = $newStatus
$dede.Range("C9").Value = $newStatus
